# Apply scheduled-runner price/profit updates to the Chocobo Profits workbook
# Each sheet corresponds to a crafting class (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Column legend: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#                K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

# ALC sheet, row 6
$ws = $wb.Worksheets.Item(1)
$ws.Range("H6").Value = 3120.2
$ws.Range("I6").Value = 134.66667
$ws.Range("J6").Value = 3866.5833
$ws.Range("K6").Value = 404.00001
$ws.Range("L6").Value = 11599.7499
$ws.Range("M6").Value = -292.00001
$ws.Range("N6").Value = -11823.7499

# ALC sheet, row 87
$ws = $wb.Worksheets.Item(1)
$ws.Range("H87").Value = 25250.666
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 25250.666
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 25250.666
$ws.Range("N87").Value = -27746.666

# ALC sheet, row 90
$ws = $wb.Worksheets.Item(1)
$ws.Range("H90").Value = 25250.666
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 25250.666
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 75751.99800000001
$ws.Range("N90").Value = -88231.99800000001

# ALC sheet, row 93
$ws = $wb.Worksheets.Item(1)
$ws.Range("H93").Value = 23159.098
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 23159.098
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 23159.098
$ws.Range("N93").Value = -28151.098

# ALC sheet, row 129
$ws = $wb.Worksheets.Item(1)
$ws.Range("H129").Value = 857.17
$ws.Range("I129").Value = 495
$ws.Range("J129").Value = 864.5612
$ws.Range("K129").Value = 1485
$ws.Range("L129").Value = 2593.6836
$ws.Range("M129").Value = 3515
$ws.Range("N129").Value = -12593.6836

# ALC sheet, row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 27302090
$ws.Range("I132").Value = 33338260
$ws.Range("J132").Value = 1432788.1
$ws.Range("K132").Value = 100014780
$ws.Range("L132").Value = 4298364.300000001
$ws.Range("M132").Value = -100012250
$ws.Range("N132").Value = -4303424.300000001

# ARM sheet, row 25
$ws = $wb.Worksheets.Item(2)
$ws.Range("H25").Value = 4173.5
$ws.Range("I25").Value = 1485.8
$ws.Range("J25").Value = 8653
$ws.Range("K25").Value = 1485.8
$ws.Range("L25").Value = 8653
$ws.Range("M25").Value = -1083.8
$ws.Range("N25").Value = -9457

# ARM sheet, row 35
$ws = $wb.Worksheets.Item(2)
$ws.Range("H35").Value = 14324.5
$ws.Range("I35").Value = 1650
$ws.Range("J35").Value = 26999
$ws.Range("K35").Value = 1650
$ws.Range("L35").Value = 26999
$ws.Range("M35").Value = -1244
$ws.Range("N35").Value = -27811

# ARM sheet, row 121
$ws = $wb.Worksheets.Item(2)
$ws.Range("H121").Value = 28214.96
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 28214.96
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 28214.96
$ws.Range("N121").Value = -31708.96

# ARM sheet, row 123
$ws = $wb.Worksheets.Item(2)
$ws.Range("H123").Value = 56399.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 56399.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 56399.5
$ws.Range("N123").Value = -66199.5

# BSM sheet, row 7
$ws = $wb.Worksheets.Item(3)
$ws.Range("H7").Value = 3656701.2
$ws.Range("I7").Value = 13745
$ws.Range("J7").Value = 4466247
$ws.Range("K7").Value = 13745
$ws.Range("L7").Value = 4466247
$ws.Range("M7").Value = -13632
$ws.Range("N7").Value = -4466473

# BSM sheet, row 11
$ws = $wb.Worksheets.Item(3)
$ws.Range("H11").Value = 1682
$ws.Range("I11").Value = 864.6667
$ws.Range("J11").Value = 1835.25
$ws.Range("K11").Value = 864.6667
$ws.Range("L11").Value = 1835.25
$ws.Range("M11").Value = -724.6667
$ws.Range("N11").Value = -2115.25

# BSM sheet, row 37
$ws = $wb.Worksheets.Item(3)
$ws.Range("H37").Value = 3009.5557
$ws.Range("I37").Value = 155.28572
$ws.Range("J37").Value = 12999.5
$ws.Range("K37").Value = 155.28572
$ws.Range("L37").Value = 12999.5
$ws.Range("M37").Value = -18.28572
$ws.Range("N37").Value = -13273.5

# BSM sheet, row 86
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 2337.3635
$ws.Range("I86").Value = 2151.375
$ws.Range("J86").Value = 2833.3333
$ws.Range("K86").Value = 2151.375
$ws.Range("L86").Value = 2833.3333
$ws.Range("M86").Value = -1028.375
$ws.Range("N86").Value = -5079.3333

# BSM sheet, row 89
$ws = $wb.Worksheets.Item(3)
$ws.Range("H89").Value = 2337.3635
$ws.Range("I89").Value = 2151.375
$ws.Range("J89").Value = 2833.3333
$ws.Range("K89").Value = 10756.875
$ws.Range("L89").Value = 14166.6665
$ws.Range("M89").Value = -5140.875
$ws.Range("N89").Value = -25398.6665

# BSM sheet, row 134
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 3073.7222
$ws.Range("I134").Value = 1620.8667
$ws.Range("J134").Value = 10338
$ws.Range("K134").Value = 4862.6001
$ws.Range("L134").Value = 31014
$ws.Range("M134").Value = -2327.6001
$ws.Range("N134").Value = -36084

# CRP sheet, row 12
$ws = $wb.Worksheets.Item(4)
$ws.Range("H12").Value = 6347.846
$ws.Range("I12").Value = 2225
$ws.Range("J12").Value = 6691.4165
$ws.Range("K12").Value = 2225
$ws.Range("L12").Value = 6691.4165
$ws.Range("M12").Value = -2055
$ws.Range("N12").Value = -7031.4165

# CRP sheet, row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2947.7708
$ws.Range("I31").Value = 1675.2439
$ws.Range("J31").Value = 10401.143
$ws.Range("K31").Value = 1675.2439
$ws.Range("L31").Value = 10401.143
$ws.Range("M31").Value = -1380.2439
$ws.Range("N31").Value = -10991.143

# CRP sheet, row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 2947.7708
$ws.Range("I34").Value = 1675.2439
$ws.Range("J34").Value = 10401.143
$ws.Range("K34").Value = 1675.2439
$ws.Range("L34").Value = 10401.143
$ws.Range("M34").Value = -1473.2439
$ws.Range("N34").Value = -10805.143

# CRP sheet, row 137
$ws = $wb.Worksheets.Item(4)
$ws.Range("H137").Value = 45286.668
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 45286.668
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 45286.668
$ws.Range("N137").Value = -55486.668

# CUL sheet, row 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 581554.6
$ws.Range("I5").Value = 409.92856
$ws.Range("J5").Value = 1485557.5
$ws.Range("K5").Value = 1229.78568
$ws.Range("L5").Value = 4456672.5
$ws.Range("M5").Value = -1117.78568
$ws.Range("N5").Value = -4456896.5

# CUL sheet, row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 775.39
$ws.Range("I131").Value = 305
$ws.Range("J131").Value = 816.29346
$ws.Range("K131").Value = 915
$ws.Range("L131").Value = 2448.88038
$ws.Range("M131").Value = 4125
$ws.Range("N131").Value = -12528.88038

# CUL sheet, row 135
$ws = $wb.Worksheets.Item(5)
$ws.Range("H135").Value = 581554.6
$ws.Range("I135").Value = 409.92856
$ws.Range("J135").Value = 1485557.5
$ws.Range("K135").Value = 3689.35704
$ws.Range("L135").Value = 13370017.5
$ws.Range("M135").Value = -1154.35704
$ws.Range("N135").Value = -13375087.5

# GSM sheet, row 70
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 6145.3076
$ws.Range("I70").Value = 5685.6
$ws.Range("J70").Value = 7677.6665
$ws.Range("K70").Value = 5685.6
$ws.Range("L70").Value = 7677.6665
$ws.Range("M70").Value = -5415.6
$ws.Range("N70").Value = -8217.666499999999

# GSM sheet, row 73
$ws = $wb.Worksheets.Item(6)
$ws.Range("H73").Value = 6145.3076
$ws.Range("I73").Value = 5685.6
$ws.Range("J73").Value = 7677.6665
$ws.Range("K73").Value = 5685.6
$ws.Range("L73").Value = 7677.6665
$ws.Range("M73").Value = -4749.6
$ws.Range("N73").Value = -9549.666499999999

# GSM sheet, row 122
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 4020.3
$ws.Range("I122").Value = 3393.2
$ws.Range("J122").Value = 5901.6
$ws.Range("K122").Value = 10179.6
$ws.Range("L122").Value = 17704.8
$ws.Range("M122").Value = -7729.599999999999
$ws.Range("N122").Value = -22604.8

# GSM sheet, row 123
$ws = $wb.Worksheets.Item(6)
$ws.Range("H123").Value = 10210.75
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10210.75
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 10210.75
$ws.Range("N123").Value = -15110.75

# GSM sheet, row 126
$ws = $wb.Worksheets.Item(6)
$ws.Range("H126").Value = 3362.073
$ws.Range("I126").Value = 2878.3098
$ws.Range("J126").Value = 4735.96
$ws.Range("K126").Value = 8634.929400000001
$ws.Range("L126").Value = 14207.88
$ws.Range("M126").Value = -6164.929400000001
$ws.Range("N126").Value = -19147.88

# GSM sheet, row 133
$ws = $wb.Worksheets.Item(6)
$ws.Range("H133").Value = 57085
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 57085
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 57085
$ws.Range("N133").Value = -67205

# LTW sheet, row 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 3390.8147
$ws.Range("I7").Value = 1762.6
$ws.Range("J7").Value = 8042.857
$ws.Range("K7").Value = 1762.6
$ws.Range("L7").Value = 8042.857
$ws.Range("M7").Value = -1650.6
$ws.Range("N7").Value = -8266.857

# LTW sheet, row 100
$ws = $wb.Worksheets.Item(7)
$ws.Range("H100").Value = 2489.2727
$ws.Range("I100").Value = 2096.1667
$ws.Range("J100").Value = 2961
$ws.Range("K100").Value = 2096.1667
$ws.Range("L100").Value = 2961
$ws.Range("M100").Value = -1555.1667
$ws.Range("N100").Value = -4043

# LTW sheet, row 122
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 6236.727
$ws.Range("I122").Value = 5067.1113
$ws.Range("J122").Value = 11500
$ws.Range("K122").Value = 15201.3339
$ws.Range("L122").Value = 34500
$ws.Range("M122").Value = -12751.3339
$ws.Range("N122").Value = -39400

# LTW sheet, row 126
$ws = $wb.Worksheets.Item(7)
$ws.Range("H126").Value = 3390.8147
$ws.Range("I126").Value = 1762.6
$ws.Range("J126").Value = 8042.857
$ws.Range("K126").Value = 5287.799999999999
$ws.Range("L126").Value = 24128.571
$ws.Range("M126").Value = -2817.799999999999
$ws.Range("N126").Value = -29068.571

# LTW sheet, row 132
$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value = 6039.6294
$ws.Range("I132").Value = 2168.9412
$ws.Range("J132").Value = 12619.8
$ws.Range("K132").Value = 6506.823600000001
$ws.Range("L132").Value = 37859.39999999999
$ws.Range("M132").Value = -3976.823600000001
$ws.Range("N132").Value = -42919.39999999999

# WVR sheet, row 47
$ws = $wb.Worksheets.Item(8)
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# WVR sheet, row 123
$ws = $wb.Worksheets.Item(8)
$ws.Range("H123").Value = 36188.75
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 36188.75
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 36188.75
$ws.Range("N123").Value = -45988.75

# WVR sheet, row 130
$ws = $wb.Worksheets.Item(8)
$ws.Range("H130").Value = 62300
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 62300
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 62300
$ws.Range("N130").Value = -72340
